# Adding Part 2 and 4 (#5)
#
# Revision-history table (Table 1) in the document has columns:
#   1) Date   2) Revision   3) Description   4) Author
#
# Row 8 (Date "9/26/2024", Revision "1.5") currently has empty
# Description and Author cells:
#   - Description -> fill it in with "Added part 2 and 4".
#   - Author      -> fill it in with "Rhenjiro Gunawan" (flagged by
#     Word's spell checker as two not-in-dictionary words, hence the
#     spellStart/spellEnd proofErr markers around each word).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Cell 1: Description for the "9/26/2024" / "1.5" row ---
$descCell = $t.Cell(8, 3)
$descCell.Range.Text = "Added part 2 and 4"

# --- Cell 2: Author for that same row ---
$authorCell = $t.Cell(8, 4)
$authorRange = $authorCell.Range
# Drop the trailing cell-mark/paragraph-mark character from the range
# so the inserted XML replaces the paragraph content in place.
$authorRange.End = $authorRange.End - 1

$authorParagraphXml = '<w:p ' + `
  'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
  'w14:paraId="18DE8D22" w14:textId="77777777" w:rsidR="008669A5" w:rsidRDefault="008669A5">' + `
    '<w:pPr>' + `
      '<w:widowControl w:val="0"/>' + `
      '<w:spacing w:line="240" w:lineRule="auto"/>' + `
    '</w:pPr>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Rhenjiro</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Gunawan</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'

$authorRange.InsertXML($authorParagraphXml)
